# Mystic Spice Premium Chai Tea product description - text content update
# (label renames + revised body copy), matching the target revision.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Para 1 label: 製品説明 -> 製品の説明
Replace-Text "製品説明" "製品の説明"

# Para 4: label + body
Replace-Text "本格ブレンド" "本格的なブレンド"
Replace-Text `
  ": 当社のチャイは、高級な紅茶葉と、シナモン、カルダモン、クローブ、ジンジャー、ブラック ペッパーなどの代表的な挽いたスパイスの調和のとれたミックスです。この古くから伝わるレシピは、一口飲むごとに本格的でしっかりとした味わいを約束します。" `
  ": プレミアムな紅茶の葉、シナモン、カルダモン、クローブ、ショウガ、黒コショウなどの厳選されたスパイスを粉砕し、絶妙にブレンドしたチャイです。この古くから伝わるレシピは、一口飲むごとに本格的でしっかりとした味わいを約束します。"

# Para 5: label + body
Replace-Text "健康増進成分" "健康に良い素材"
Replace-Text `
  ": Mystic Spice Chai Tea の各成分は、自然な健康上の利点を考慮して選択されています。ジンジャーとカルダモンは消化を助け、シナモンは血糖値の調節を助け、クローブは抗酸化物質を高めます。" `
  ": Mystic Spice Chai Tea には、健康効果を考慮して選ばれた天然素材が使用されています。ジンジャーとカルダモンは消化を助け、シナモンは血糖値の調節を助け、クローブは抗酸化物質を高めます。"

# Para 6: body only (label 豊かな香りと風味 unchanged)
Replace-Text `
  ": 温かくスパイシーな香りと深く爽快な味わいのチャイは、一日の始まりや夜のくつろぎに最適な飲み物です。風味は強烈でありながらバランスが取れており、快適で心地よい体験を生み出します。" `
  ": 温かくスパイシーな香りと深く爽快な味わいで、一日の始まりや夜のリラックスタイムに最適な飲み物です。風味は強烈でありながらバランスが取れており、快適で心地よい体験を生み出します。"

# Para 7: label + body
Replace-Text "多彩な淹れ方" "多様な楽しみ方"
Replace-Text `
  ": 熱々のチャイが好きでも、さわやかなアイス ティーとしても、クリーミーなラテとしても、当社のブレンドはどんな好みにも合う多用途な製品です。お好みの方法でチャイをお楽しみいただけるよう、簡単な淹れ方の説明書が付属しています。" `
  ": ホットだけでなく、爽やかなアイスティーやクリーミーなラテなど、好みに合わせて自由にお楽しみいただける万能なブレンドです。お好みの方法でチャイをお楽しみいただけるよう、簡単な淹れ方の説明書が付属しています。"

# Para 8: body only (label 持続可能な調達 unchanged)
Replace-Text `
  ": 持続可能性を重視し、有機農業を実践する小規模農場から原材料を調達し、最高の品質だけでなく地球の福祉も保証します。" `
  ": 持続可能性を重視し、有機農業を実践する小規模農場から原料を調達することで、最高品質を実現するだけでなく地球環境にも配慮しています。"

# Para 9: body only (label エレガントなパッケージ unchanged)
Replace-Text `
  ": Mystic Spice Chai Tea は、美しくデザインされた環境に優しいパッケージに入っており、紅茶愛好家への贈り物や自分への贅沢なご褒美に最適です。" `
  ": Mystic Spice Chai Tea は、環境に配慮した美しいデザインのパッケージに梱包されているため、お茶が好きな方に贈るギフトや自分自身への贅沢なご褒美として最適です。"

# Para 10: body only (label 顧客満足度保証 unchanged)
Replace-Text `
  ": 当社は自社製品に責任を持ち、顧客満足度を保証します。Mystic Spice Chai Tea がお客様のご期待に添えない場合は、当社が改善するよう努めます。" `
  ": 当社は製品の品質に自信を持っており、満足度保証を提供しています。Mystic Spice Chai Tea がお客様のご期待に添えない場合は、当社が改善するよう努めます。"

# Para 11: label + body
Replace-Text "次のお客様に最適です" "本製品が最適な方"
Replace-Text `
  ": 紅茶愛好家、健康志向の人、温かくてスパイシーな飲み物の愛好家、そして伝統的なインドのチャイの豊かな風味を探求したいお客様など。" `
  ": お茶が好きな方、健康志向の方、温かくスパイシーな飲み物が好きな方、伝統的なインドのチャイの豊かな風味を体験したい方。"

# Para 12: closing paragraph
Replace-Text `
  "Mystic Spice Premium Chai Tea で本場のインドの味をお楽しみください。すべての一杯が風味と伝統の物語です。" `
  "Mystic Spice Premium Chai Tea でインドの本格的な風味をお楽しみください。一杯ごとに豊かな風味と伝統を感じられます。"

Write-Host "Done"
